# Team_Everyday_Attendence.xlsx update
# - Append attendance rows for 24-Aug-2023 (row 24) and 25-Aug-2023 (row 25)
# - Add reviewer comments explaining the new ABSENT marks
# - Move the active selection down to the newly entered data (K25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 24 : 24-Aug-2023 -------------------------------------------------
$ws.Range("A24").Value = 45162
$ws.Range("A24").NumberFormat = "d-mmm-yy"
$ws.Range("B24").Value = "PRESENT"
$ws.Range("C24").Value = "PRESENT"
$ws.Range("D24").Value = "PRESENT"
$ws.Range("E24").Value = "PRESENT"
$ws.Range("F24").Value = "PRESENT"
$ws.Range("G24").Value = "PRESENT"
$ws.Range("H24").Value = "PRESENT"
$ws.Range("I24").Value = "PRESENT"
$ws.Range("J24").Value = "ABSENT"
$ws.Range("K24").Value = "PRESENT"

# ---- Row 25 : 25-Aug-2023 -------------------------------------------------
$ws.Range("A25").Value = 45163
$ws.Range("A25").NumberFormat = "d-mmm-yy"
$ws.Range("B25").Value = "PRESENT"
$ws.Range("C25").Value = "PRESENT"
$ws.Range("D25").Value = "PRESENT"
$ws.Range("E25").Value = "PRESENT"
$ws.Range("F25").Value = "PRESENT"
$ws.Range("G25").Value = "ABSENT"
$ws.Range("H25").Value = "PRESENT"
$ws.Range("I25").Value = "PRESENT"
$ws.Range("J25").Value = "ABSENT"
$ws.Range("K25").Value = "PRESENT"

# ---- Comments explaining the ABSENT marks --------------------------------
$ws.Range("J24").AddComment("LENOVO:" + [char]10 + "No response")
$ws.Range("G25").AddComment("LENOVO:" + [char]10 + "No Response")
$ws.Range("J25").AddComment("LENOVO:" + [char]10 + "No Response")

# ---- Move selection to the last entered cell ------------------------------
$ws.Range("K25").Select()
